$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Expand the parts table with two more data rows -----------------------
# Row 9 is currently a blank data row (already inside the table) that needs
# to be filled in. Two brand-new rows need to be inserted ahead of the
# totals row (currently row 10) so the totals row ends up at row 12.
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G12"))

# --- Row 9: CPC - DC power receptacle connector ----------------------------
# (Item text is entered before Seller so the new shared strings land in the
# same order as the authored workbook: CLIFF... then CPC.)
$ws.Range("B9").Value = "CLIFF ELECTRONIC COMPONENTS - FC681473 - CONNECTOR, RECEPTACLE, DC POWER, 2.1MM"
$ws.Range("A9").Value = "CPC"
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0.38
$ws.Range("E9").Formula = "=C9*D9"
$ws.Range("F9").Value = 1.88
$ws.Range("G9").Formula = "=E9+F9"

# --- Row 10: eBay - Strut hinges --------------------------------------------
$ws.Range("A10").Value = "eBay"
$ws.Range("B10").Value = "Strut Hinges For Flight Case Lids - One Pair - Nickel Finish"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 4.49
$ws.Range("E10").Formula = "=C10*D10"
$ws.Range("F10").Value = $null
$ws.Range("G10").Formula = "=E10+F10"

# --- Row 11: eBay - Clear acrylic sheet ------------------------------------
$ws.Range("A11").Value = "eBay"
$ws.Range("B11").Value = "Clear acrylic 4mmx100mmx200mm"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 3.16
$ws.Range("E11").Formula = "=C11*D11"
$ws.Range("F11").Value = $null
$ws.Range("G11").Formula = "=E11+F11"

# Match cell formatting used by the other eBay rows (e.g. row 8).
$ws.Range("A8:G8").Copy()
$ws.Range("A10:G10").PasteSpecial(-4122)
$ws.Range("A11:G11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the selection Excel ends up with after this kind of edit.
$ws.Range("E18").Select()

$wb.Save()
